$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '22.425.82'
$ws.Range("E2").Value = '  +8.76%  '
$ws.Range("D3").Value = '1.602.14'
$ws.Range("E3").Value = '  +8.11%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.004'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.51%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.9913'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.14%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '300.57'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +7.51%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3686'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.52%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3412'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +10.65%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '42.56'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +6.17%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.141'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +7.02%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07055'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +5.49%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.9999'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.39%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.945'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +7.37%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '19.73'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +8.89%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.626'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +6.48%  '
$ws.Range("E16").Value = '  +5.12%  '
$ws.Range("D17").Value = '1.597.87'
$ws.Range("E17").Value = '  +8.10%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.9916'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.13%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06766'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +13.61%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '77.94'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +11.71%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '16.13'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +10.75%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.029'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +9.27%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.81'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +6.54%  '
$ws.Range("D24").Value = '22.442.99'
$ws.Range("E24").Value = '  +8.64%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.403'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +5.88%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.529'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +17.97%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '150.71'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +5.97%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.56'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +12.94%  '
$ws.Range("D29").Value = '1.780.50'
$ws.Range("E29").Value = '  +8.75%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '121.11'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +5.78%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.187'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +6.53%  '
$ws.Range("E32").Value = '  +20.07%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.9522'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +15.28%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08271'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.42%  '
$ws.Range("E35").Value = '  +6.30%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.302'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +11.21%  '
$ws.Range("B37").Value = 'Aptos'
$ws.Range("C37").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '11.95'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +13.81%  '
$ws.Range("B38").Value = 'TrustWalletToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.268'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +4.88%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.566'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +12.14%  '
$ws.Range("E40").Value = '  +5.33%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.02217'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +8.02%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.2027'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +7.58%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.9917'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.19%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.5921'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +11.27%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.800'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +7.36%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '13.23'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +8.24%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5691'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +9.16%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '127.38'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +7.31%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.971'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +8.70%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06814'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +4.86%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '73.89'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +8.94%  '
